# Update attendance / price figures and refresh the "2024-12-27" event dates
# to "2025-01-01" (now that the event has happened again / been rescheduled)
# across the relevant sheets of the workbook.
#
# Sheet layout: 1=展览, 2=演出, 3=本地生活 (empty), 4=全部类型

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value  = 15778
$ws1.Range("F8").Value  = 709
$ws1.Range("F9").Value  = 15463
$ws1.Range("F11").Value = 9067
$ws1.Range("F16").Value = 204
$ws1.Range("F21").Value = 559
$ws1.Range("F25").Value = 1117
$ws1.Range("F33").Value = 64
$ws1.Range("F34").Value = 45
$ws1.Range("F35").Value = 260
$ws1.Range("F37").Value = 461
$ws1.Range("F39").Value = 5571

# Row 40: event re-dated from 2024-12-27 to 2025-01-01, and it is now sellable
# (text "不可售" -> numeric min price 68). Force B40 to stay plain text (it
# would otherwise be auto-parsed into a date serial by Excel).
$ws1.Range("B40").Value = "'2025-01-01"
$ws1.Range("E40").Value = "2025.01.01 09:00-01.02 16:00"
$ws1.Range("G40").Value = 68

# ---------------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 72

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 15778
$ws4.Range("F8").Value  = 709
$ws4.Range("F9").Value  = 15463
$ws4.Range("F11").Value = 9068
$ws4.Range("F16").Value = 204
$ws4.Range("F21").Value = 559
$ws4.Range("F25").Value = 1117
$ws4.Range("F32").Value = 72
$ws4.Range("F35").Value = 64
$ws4.Range("F36").Value = 45
$ws4.Range("F37").Value = 260
$ws4.Range("F39").Value = 461
$ws4.Range("F41").Value = 5571

# Row 43: same re-date / re-price as row 40 on the "展览" sheet above.
$ws4.Range("B43").Value = "'2025-01-01"
$ws4.Range("E43").Value = "2025.01.01 09:00-01.02 16:00"
$ws4.Range("G43").Value = 68
